$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

# Update D3:D5 status text from "70% (Xong 2/3 Chức năng)" to "70% (Chưa có cn Cập Nhật)"
$ws.Range("D3").Value = "70% (Chưa có cn Cập Nhật)"
$ws.Range("D4").Value = "70% (Chưa có cn Cập Nhật)"
$ws.Range("D5").Value = "70% (Chưa có cn Cập Nhật)"

# Hoang integrated QuanLyLoaiHang function at 100% - set final result for row 6
$ws.Range("E6").Value = "100% (9/06/2010)"
